$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 96 ("Rental and leasing services" / 532),
# shifting all subsequent rows down by one.
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row with the new NAICS crosswalk entry.
$ws.Cells.Item(96, 1).Value = "Other real estate activities"
$ws.Cells.Item(96, 2).Value = "5313"

# Match the author's saved view state: scrolled down with B97 selected.
$ws.Activate()
$ws.Range("B97").Select()
